$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.339.95'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').Value = '1.846.35'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9977'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.07'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6264'
$ws.Range('E6').Value = '  -0.64%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9986'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07600'
$ws.Range('E8').Value = '  -1.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2900'
$ws.Range('E9').Value = '  -1.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.72'
$ws.Range('E10').Value = '  +0.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07731'
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.00001061'
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '82.91'
$ws.Range('E15').Value = '  -1.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.131'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('D17').Value = '29.372.48'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '227.43'
$ws.Range('E18').Value = '  -0.86%  '
$ws.Range('E19').Value = '  -1.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9986'
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '158.50'
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('E24').Value = '  -0.53%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.426'
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.64'
$ws.Range('E26').Value = '  -0.28%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.419'
$ws.Range('E27').Value = '  +7.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.459'
$ws.Range('E28').Value = '  -0.65%  '
$ws.Range('E29').Value = '  -2.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.101'
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.161'
$ws.Range('E32').Value = '  +0.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.831'
$ws.Range('E33').Value = '  -1.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.6948'
$ws.Range('E34').Value = '  -2.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.582'
$ws.Range('E35').Value = '  -0.24%  '
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('D37').Value = '1.226.40'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.718'
$ws.Range('E38').Value = '  -2.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.352'
$ws.Range('E39').Value = '  -1.73%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8955'
$ws.Range('E40').Value = '  -1.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9983'
$ws.Range('E42').Value = '  -0.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '65.45'
$ws.Range('E43').Value = '  -1.34%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.185'
$ws.Range('E44').Value = '  +0.32%  '
$ws.Range('E45').Value = '  -0.99%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.688'
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.985'
$ws.Range('E47').Value = '  -0.67%  '
$ws.Range('E48').Value = '  -4.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1140'
$ws.Range('E49').Value = '  +1.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05697'
$ws.Range('E50').Value = '  -0.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4619'
$ws.Range('E51').Value = '  -0.18%  '
